# Apply "generación actividad 1 json" data refresh to cleaned_data_name sheet.
# Updates marketCapUsd (F) for every row, plus supply (D) figures that changed,
# and re-sorts a handful of neighbouring rows whose market-cap ordering flipped
# (id/name/supply/explorer swap together with the refreshed marketCapUsd).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row = 2; Col = 6; Value = 1646598735231.417},
    @{Row = 3; Col = 6; Value = 216950970985.5372},
    @{Row = 4; Col = 6; Value = 144167887216.4139},
    @{Row = 5; Col = 6; Value = 120175935170.5374},
    @{Row = 6; Col = 6; Value = 85330557169.34145},
    @{Row = 7; Col = 6; Value = 60828634979.1779},
    @{Row = 8; Col = 6; Value = 59533179397.05877},
    @{Row = 9; Col = 6; Value = 24087585742.23228},
    @{Row = 10; Col = 6; Value = 22907170408.59418},
    @{Row = 11; Col = 6; Value = 22706495702.73655},
    @{Row = 12; Col = 4; Value = 9374253.737797299},
    @{Row = 12; Col = 6; Value = 16837482100.90499},
    @{Row = 13; Col = 6; Value = 10657415753.52441},
    @{Row = 14; Col = 6; Value = 8800913359.149023},
    @{Row = 15; Col = 6; Value = 8721616735.176374},
    @{Row = 16; Col = 6; Value = 8403976328.942166},
    @{Row = 17; Col = 6; Value = 8010382263.904924},
    @{Row = 18; Col = 6; Value = 7513765081.376249},
    @{Row = 19; Col = 6; Value = 7195598533.846197},
    @{Row = 20; Col = 6; Value = 7159684628.490468},
    @{Row = 21; Col = 6; Value = 6774138057.537104},
    @{Row = 22; Col = 6; Value = 6297023774.210812},
    @{Row = 23; Col = 6; Value = 6191043531.272076},
    @{Row = 24; Col = 6; Value = 6129985585.660106},
    @{Row = 25; Col = 6; Value = 6085601183.095393},
    @{Row = 26; Col = 6; Value = 5987629528.959108},
    @{Row = 27; Col = 6; Value = 5400664667.563888},
    @{Row = 28; Col = 6; Value = 5381097135.256779},
    @{Row = 29; Col = 6; Value = 5364279939.315865},
    @{Row = 30; Col = 6; Value = 5248524751.149027},
    @{Row = 31; Col = 1; Value = 'monero'},
    @{Row = 31; Col = 2; Value = 'monero'},
    @{Row = 31; Col = 4; Value = 18446744.07370955},
    @{Row = 31; Col = 5; Value = 'http://moneroblocks.info/'},
    @{Row = 31; Col = 6; Value = 3952889456.807546},
    @{Row = 32; Col = 1; Value = 'hyperliquid'},
    @{Row = 32; Col = 2; Value = 'hyperliquid'},
    @{Row = 32; Col = 4; Value = 333928180},
    @{Row = 32; Col = 5; Value = 'https://app.hyperliquid.xyz/explorer/token/0x0d01dc56dcaaca66ad901c959b4011ec'},
    @{Row = 32; Col = 6; Value = 3942293145.254296},
    @{Row = 33; Col = 6; Value = 3802060652.464745},
    @{Row = 34; Col = 6; Value = 3681192143.871577},
    @{Row = 35; Col = 6; Value = 3036171573.969729},
    @{Row = 36; Col = 6; Value = 2987845446.324038},
    @{Row = 37; Col = 6; Value = 2833157119.960022},
    @{Row = 38; Col = 6; Value = 2811123801.963756},
    @{Row = 39; Col = 6; Value = 2586208908.081504},
    @{Row = 40; Col = 1; Value = 'mantle'},
    @{Row = 40; Col = 2; Value = 'mantle'},
    @{Row = 40; Col = 4; Value = 3364694382.836841},
    @{Row = 40; Col = 5; Value = 'https://mantlescan.xyz/'},
    @{Row = 40; Col = 6; Value = 2509934872.143909},
    @{Row = 41; Col = 1; Value = 'crypto-com-coin'},
    @{Row = 41; Col = 2; Value = 'cronos'},
    @{Row = 41; Col = 4; Value = 26571560696},
    @{Row = 41; Col = 5; Value = 'https://etherscan.io/token/0xa0b73e1ff0b80914ab6fe0444e65848c4c34450b'},
    @{Row = 41; Col = 6; Value = 2508854481.305764},
    @{Row = 42; Col = 6; Value = 2478927769.463691},
    @{Row = 43; Col = 6; Value = 2444652886.993679},
    @{Row = 44; Col = 6; Value = 2351184089.48359},
    @{Row = 45; Col = 6; Value = 2254443603.579769},
    @{Row = 46; Col = 6; Value = 1978702466.497413},
    @{Row = 47; Col = 6; Value = 1895039081.183438},
    @{Row = 48; Col = 6; Value = 1886540877.231068},
    @{Row = 49; Col = 6; Value = 1866661757.112612},
    @{Row = 50; Col = 6; Value = 1846890920.009128},
    @{Row = 51; Col = 6; Value = 1840368304.463784},
    @{Row = 52; Col = 6; Value = 1761381396.987541},
    @{Row = 53; Col = 6; Value = 1698217304.890254},
    @{Row = 54; Col = 6; Value = 1676857733.536231},
    @{Row = 55; Col = 6; Value = 1665964719.856101},
    @{Row = 56; Col = 6; Value = 1655783569.038651},
    @{Row = 57; Col = 6; Value = 1540096770.409258},
    @{Row = 58; Col = 6; Value = 1433271486.071678},
    @{Row = 59; Col = 6; Value = 1386990752.048404},
    @{Row = 60; Col = 6; Value = 1313792591.457205},
    @{Row = 61; Col = 6; Value = 1286879157.824451},
    @{Row = 62; Col = 6; Value = 1246753450.860751},
    @{Row = 63; Col = 6; Value = 1137944283.31947},
    @{Row = 64; Col = 6; Value = 1094048852.565082},
    @{Row = 65; Col = 6; Value = 1091583832.133298},
    @{Row = 66; Col = 6; Value = 1088034541.022287},
    @{Row = 67; Col = 6; Value = 1051758568.414132},
    @{Row = 68; Col = 6; Value = 1042373555.02302},
    @{Row = 69; Col = 6; Value = 957543131.8998387},
    @{Row = 70; Col = 6; Value = 952334860.657215},
    @{Row = 71; Col = 6; Value = 931258904.8206109},
    @{Row = 72; Col = 6; Value = 882723897.8018025},
    @{Row = 73; Col = 6; Value = 850694144.3058498},
    @{Row = 74; Col = 6; Value = 843406096.0131915},
    @{Row = 75; Col = 1; Value = 'injective-protocol'},
    @{Row = 75; Col = 2; Value = 'injective'},
    @{Row = 75; Col = 4; Value = 98970935.41},
    @{Row = 75; Col = 5; Value = 'https://etherscan.io/token/0xe28b3b32b6c345a34ff64674606124dd5aceca30'},
    @{Row = 75; Col = 6; Value = 830844404.7240243},
    @{Row = 76; Col = 1; Value = 'binaryx-new'},
    @{Row = 76; Col = 2; Value = 'four'},
    @{Row = 76; Col = 4; Value = 381867255.144574},
    @{Row = 76; Col = 5; Value = 'Sin datos'},
    @{Row = 76; Col = 6; Value = 829706005.3515067},
    @{Row = 77; Col = 6; Value = 824480212.14785},
    @{Row = 78; Col = 6; Value = 820486481.2163347},
    @{Row = 79; Col = 6; Value = 805832036.9772431},
    @{Row = 80; Col = 6; Value = 805652785.9603587},
    @{Row = 81; Col = 6; Value = 794664818.6886672},
    @{Row = 82; Col = 6; Value = 781445825.702072},
    @{Row = 83; Col = 6; Value = 770838100.0799955},
    @{Row = 84; Col = 6; Value = 757701482.3358722},
    @{Row = 85; Col = 6; Value = 747448705.8206316},
    @{Row = 86; Col = 6; Value = 681698637.0396736},
    @{Row = 87; Col = 6; Value = 680596986.34711},
    @{Row = 88; Col = 6; Value = 674734612.2559657},
    @{Row = 89; Col = 6; Value = 656166222.6006479},
    @{Row = 90; Col = 1; Value = 'gala'},
    @{Row = 90; Col = 2; Value = 'gala'},
    @{Row = 90; Col = 4; Value = 43820741411.30103},
    @{Row = 90; Col = 5; Value = 'https://ethplorer.io/es/address/0x15d4c048f83bd7e37d49ea4c83a07267ec4203da#chart=candlestick'},
    @{Row = 90; Col = 6; Value = 652769598.9738292},
    @{Row = 91; Col = 1; Value = 'kaia'},
    @{Row = 91; Col = 2; Value = 'kaia'},
    @{Row = 91; Col = 4; Value = 5966455361.82},
    @{Row = 91; Col = 5; Value = 'https://www.kaiascan.io/'},
    @{Row = 91; Col = 6; Value = 651584475.9770766},
    @{Row = 92; Col = 1; Value = 'the-sandbox'},
    @{Row = 92; Col = 2; Value = 'the sandbox'},
    @{Row = 92; Col = 4; Value = 2481357126.223322},
    @{Row = 92; Col = 5; Value = 'https://etherscan.io/token/0x3845badAde8e6dFF049820680d1F14bD3903a5d0'},
    @{Row = 92; Col = 6; Value = 649486304.218258},
    @{Row = 93; Col = 1; Value = 'berachain'},
    @{Row = 93; Col = 2; Value = 'berachain'},
    @{Row = 93; Col = 4; Value = 107480000},
    @{Row = 93; Col = 5; Value = 'Sin datos'},
    @{Row = 93; Col = 6; Value = 647760253.0012007},
    @{Row = 94; Col = 6; Value = 619351409.6742076},
    @{Row = 95; Col = 6; Value = 618858681.5296142},
    @{Row = 96; Col = 6; Value = 617379057.857505},
    @{Row = 97; Col = 6; Value = 608853422.7302326},
    @{Row = 98; Col = 6; Value = 606575928.6660298},
    @{Row = 99; Col = 6; Value = 582062382.8490387},
    @{Row = 100; Col = 6; Value = 536249794.0901753},
    @{Row = 101; Col = 6; Value = 534658843.8775095}
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, $change.Col).Value = $change.Value
}
